# Exponential_Distribution.xlsx: wire up the live EXPON.DIST() formulas that
# were previously stored as plain 0-literals, and move the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: density f(A2) = EXPON.DIST(A2, 1/mean, FALSE) using the mean in $E$1.
$ws.Range("B2").Formula = '=_xlfn.EXPON.DIST(A2,1/$E$1,FALSE)'

# B3:B26: same density formula, filled down (becomes one shared formula
# group anchored at B3, si=0, like Excel's own fill-down does).
$ws.Range("B3:B26").Formula = '=_xlfn.EXPON.DIST(A3,1/$E$1,FALSE)'

# E3: P(3 < X < 7) = CDF(7) - CDF(3).
$ws.Range("E3").Formula = '=_xlfn.EXPON.DIST(7,1/E1,TRUE)-_xlfn.EXPON.DIST(3,1/E1,TRUE)'

# The sheet's last active-cell selection moved from B2 to E2.
$ws.Range("E2").Select()

$wb.Save()
